# Apio (Feria Lagunitas de Puerto Montt) - weekly update.
# A new price observation is inserted as row 79, pushing the existing
# rows 79-169 down to 80-170 (dimension grows from A1:R169 to A1:R170).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 79; this shifts rows 79..169
# down to 80..170 and carries formatting (incl. the date style on column D)
# down with them.
$ws.Rows.Item(79).Insert()

# Fill in the new record at row 79.
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 44539
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112017
$ws.Range("G79").Value = "Apio"
$ws.Range("H79").Value = "Americana (o)"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 25
$ws.Range("K79").Value = 12000
$ws.Range("L79").Value = 12000
$ws.Range("M79").Value = 12000
$ws.Range("N79").Value = "$/docena de matas"
$ws.Range("O79").Value = "Región de Coquimbo"
$ws.Range("P79").Value = 2000
$ws.Range("Q79").Value = 6
$ws.Range("R79").Value = "Hortaliza"
